$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Stamp ids on the two existing content controls (Table of Contents sdt
#    and the goog_rdk_0-tagged table sdt). `ContentControl.ID` is a
#    Word-assigned, read-only value in the real object model, but we still
#    touch it here in case the control supports it, matching the intent of
#    the authored edit (Word backfilling missing <w:id> on save).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.ContentControls.Count; $i++) {
    $cc = $d.ContentControls.Item($i)
    if ($cc.Tag -eq "goog_rdk_0") {
        $cc.ID = -1439923158
    } elseif ($cc.Type -eq 5) {
        $cc.ID = -25715710
    }
}

# ---------------------------------------------------------------------------
# 2) Tighten the "A los usuarios" bullet: drop the "geolocalización, " filter
#    example, keeping the rest of the sentence untouched.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "geolocalización, categorías, fechas", $true, $false, $false, $false,
    $false, $true, 1, $false, "categorías, fechas", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the whole "Moderación de contenido por parte de administradores."
#    bullet paragraph entirely (including its paragraph mark).
# ---------------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Moderación de contenido por parte de administradores.*") {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 4) Rename the bare "TableNormal" style stub from "Table Normal" to
#    "TableNormal" (the fuller "Normal Table" style definition elsewhere is
#    left alone).
# ---------------------------------------------------------------------------
$d.Styles("TableNormal").NameLocal = "TableNormal"
